$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 edits -----------------------------------------------------

# New session added to slot 1 (M4): Conrau-Lewis AI pedagogy talk
$ws1.Range("B4").Copy()
$ws1.Range("M4").PasteSpecial(-4122)
$ws1.Range("M4").Value = "Conrau-Lewis: Overcoming students’ initial reactions to AI through text-based experiments"

# Presenter-name swap for the "Ungrading" / "Place-Based Education" sessions
$ws1.Range("I5").Value = "Matlack: Place-Based Education in the Era of AI"
$ws1.Range("J6").Value = "Cotton: Student Buy-In and ""Ungrading"" in the Humanities Classroom"

$ws1.Range("M4").Select()

# --- Add Sheet2 (presenter directory, inserted right after Sheet1) ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Slot 1 (10:15-11:15)"
$ws2.Range("B1").Value = "Slot 2 (11:30-12:30)"
$ws2.Range("C1").Value = "Slot 3 (1:30-2:30)"

$ws2.Range("A2").Value = "Paul: Annie Murphy Paul Breakout Session: Deeper Dive into ""Thinking with our Surroundings"" from The Extended Mind"

$ws2.Range("B3").Value = "Hojinicki: From Curiosity to Community: Launching a Teacher-Led Gen AI Professional Development Initiative"
$ws2.Range("C3").Value = "Hojinicki: From Curiosity to Community: Launching a Teacher-Led Gen AI Professional Development Initiative"

$ws2.Range("A4").Value = "Nilsson: What Learning Science Tells Us about Teaching with AI"
$ws2.Range("B4").Value = "Nilsson: Leverage AI to Support Teaching Your Passions"
$ws2.Range("C4").Value = "Nilsson: How To Define and Draw Clear Lines for Appropriate AI Use"

$ws2.Range("A5").Value = "Poole: Expanding Inquiry: Using AI Tools to Create Simulations and Investigations"
$ws2.Range("B5").Value = "LaForest: Slavery and Loomis Chaffee: An Ethical History Project: Collaborating with Students to Make Meaning in and out of the Classroom"
$ws2.Range("C5").Value = "Poole: Expanding Inquiry: Using AI Tools to Create Simulations and Investigations"

$ws2.Range("A6").Value = "Alsamadisi: Rethinking Creative and Critical Thinking in the Age of AI"
$ws2.Range("B6").Value = "Sadowitz: Low-Tech Reading: Strategies for Engaging High-Tech Screenagers"
$ws2.Range("C6").Value = "Sadowitz: Low-Tech Reading: Strategies for Engaging High-Tech Screenagers"

$ws2.Range("A7").Value = "Heckman: Lessons Learned from (Co-Lab)orating Across Schools"
$ws2.Range("B7").Value = "Ellinwood: Digital Dependency vs. Human Connection: Educational Strategies for the AI Character Era"
$ws2.Range("C7").Value = "Ellinwood: Digital Dependency vs. Human Connection: Educational Strategies for the AI Character Era"

$ws2.Range("A8").Value = "Seyboth: Better Together: How Human Connection Transforms AI into an Effective Educational Partner"
$ws2.Range("B8").Value = "Sperber: Ctrl+Alt+Engage: Rebooting Classrooms with AI Simulations"
$ws2.Range("C8").Value = "Lincoln: Build-a-Bot Workshop: Make Your Own AI to Make Sense of AI"

$ws2.Range("A9").Value = "Cotton: Student Buy-In and ""Ungrading"" in the Humanities Classroom"
$ws2.Range("B9").Value = "Cotton: Student Buy-In and ""Ungrading"" in the Humanities Classroom"
$ws2.Range("C9").Value = "Chew: Pre-Skilling for the Unknown: Building a Human-Centered AI Culture in Schools"

$ws2.Range("A10").Value = "Napirowska: The Value of Struggle: Preserving Meaningful Learning in an AI Age"
$ws2.Range("B10").Value = "Napirowska: The Value of Struggle: Preserving Meaningful Learning in an AI Age"
$ws2.Range("C10").Value = "Matlack: Place-Based Education in the Era of AI"

$ws2.Range("A11").Value = "Lamb: AI-Powered Pedagogy: Crafting Effective Prompts for Transformative Learning"
$ws2.Range("B11").Value = "Lamb: AI-Powered Pedagogy: Crafting Effective Prompts for Transformative Learning"
$ws2.Range("C11").Value = "Spaletta: AI-Adapted Writing Assignments for Skill Development and AI Literacy"

$ws2.Range("A12").Value = "Solomon: Untethered Thinking: Design Thinking Frameworks for Tech-Free Student Brainstorms"
$ws2.Range("B12").Value = "Peterson: Metacognition and AI: Empowering High School Students to Become Strategic Learners in a Technological World"
$ws2.Range("C12").Value = "Peterson: Metacognition and AI: Empowering High School Students to Become Strategic Learners in a Technological World"

$ws2.Range("B13").Value = "Morrell: What's in a Name: Metaphors We Write By"

$ws2.Range("B14").Value = "Scullin: Creating a Student Learning Praxis in the Age of AI"
$ws2.Range("C14").Value = "MacClintic: Teaching Student Presentation Skills"

# --- Formatting: reuse Sheet1's existing styles (header / body / blank) ---

# Style "2" (bold centered header w/ border) -> row 1
$ws1.Range("A4").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# Style "4" (bold centered, no border, blank) -> B2:C2
$ws1.Range("B5").Copy()
$ws2.Range("B2:C2").PasteSpecial(-4122)

# Style "3" (wrap text body) -> everything else that holds/represents body content
$ws1.Range("B4").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Range("B3:C3").PasteSpecial(-4122)
$ws2.Range("A4:C12").PasteSpecial(-4122)
$ws2.Range("B13").PasteSpecial(-4122)
$ws2.Range("A14:C14").PasteSpecial(-4122)

# Row heights (matches the source row's wrapped text extent)
$ws2.Rows("2").RowHeight = 169
$ws2.Rows("3").RowHeight = 141
$ws2.Rows("4").RowHeight = 99
$ws2.Rows("5").RowHeight = 211
$ws2.Rows("6").RowHeight = 99
$ws2.Rows("7").RowHeight = 141
$ws2.Rows("8").RowHeight = 141
$ws2.Rows("9").RowHeight = 127
$ws2.Rows("10").RowHeight = 99
$ws2.Rows("11").RowHeight = 113
$ws2.Rows("12").RowHeight = 155
$ws2.Rows("13").RowHeight = 71
$ws2.Rows("14").RowHeight = 85

$ws2.Range("C3:C14").Select()

$ws1.Activate()
$ws1.Range("M4").Select()
